$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-19 (A1:M19 data index 8-17) reflect a re-run of the averaged-intensity
# calculation including new spiral sampling schemes; rows 3-9 (and the header rows) are unchanged.

$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.066167906916869
$ws.Range("D10").Value = 0.8099588618061759
$ws.Range("E10").Value = 1.027829301863287
$ws.Range("F10").Value = 1.066167906916869
$ws.Range("G10").Value = 0.8952340248181317
$ws.Range("H10").Value = 1.085230290652888
$ws.Range("I10").Value = 1.040276857054268
$ws.Range("J10").Value = 0.8099588618061759
$ws.Range("K10").Value = 0.9188940818347316
$ws.Range("L10").Value = 0.9925309943758004
$ws.Range("M10").Value = 0.9874495405186033

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9479613986047402
$ws.Range("D11").Value = 0.9430539209197296
$ws.Range("E11").Value = 1.039176103897089
$ws.Range("F11").Value = 0.9479613986047402
$ws.Range("G11").Value = 0.9349590174498704
$ws.Range("H11").Value = 1.147763790373659
$ws.Range("I11").Value = 1.013817346987879
$ws.Range("J11").Value = 0.9430539209197296
$ws.Range("K11").Value = 0.9911150124084094
$ws.Range("L11").Value = 0.9695382055065748
$ws.Range("M11").Value = 1.004455263038828

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9478454218904558
$ws.Range("D12").Value = 0.9447663093989092
$ws.Range("E12").Value = 1.038600983981753
$ws.Range("F12").Value = 0.9478454218904558
$ws.Range("G12").Value = 0.9360565720399935
$ws.Range("H12").Value = 1.146066676967094
$ws.Range("I12").Value = 1.013301727925444
$ws.Range("J12").Value = 0.9447663093989092
$ws.Range("K12").Value = 0.9916836466903309
$ws.Range("L12").Value = 0.9697645342903934
$ws.Range("M12").Value = 1.004439615367275

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9479067629627562
$ws.Range("D13").Value = 0.9434117273001003
$ws.Range("E13").Value = 1.039070717087051
$ws.Range("F13").Value = 0.9479067629627562
$ws.Range("G13").Value = 0.9352068345844746
$ws.Range("H13").Value = 1.147421995125292
$ws.Range("I13").Value = 1.013663452452542
$ws.Range("J13").Value = 0.9434117273001003
$ws.Range("K13").Value = 0.9912412221935754
$ws.Range("L13").Value = 0.9695739925781658
$ws.Range("M13").Value = 1.004446914918703

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9646159999999996
$ws.Range("D14").Value = 0.7143520000000005
$ws.Range("E14").Value = 1.115275999999997
$ws.Range("F14").Value = 0.9646159999999996
$ws.Range("G14").Value = 0.7909080000000023
$ws.Range("H14").Value = 1.367271999999999
$ws.Range("I14").Value = 1.087151999999998
$ws.Range("J14").Value = 0.7143520000000005
$ws.Range("K14").Value = 0.9148139999999987
$ws.Range("L14").Value = 0.9397149999999992
$ws.Range("M14").Value = 1.006595999999999

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.99
$ws.Range("D15").Value = 0.37
$ws.Range("E15").Value = 1.23
$ws.Range("F15").Value = 0.99
$ws.Range("G15").Value = 0.58
$ws.Range("H15").Value = 1.69
$ws.Range("I15").Value = 1.19
$ws.Range("J15").Value = 0.37
$ws.Range("K15").Value = 0.8
$ws.Range("L15").Value = 0.895
$ws.Range("M15").Value = 1.008333333333333

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9949522190336001
$ws.Range("D16").Value = 0.6299189688320024
$ws.Range("E16").Value = 1.130944102400001
$ws.Range("F16").Value = 0.9949522190336001
$ws.Range("G16").Value = 0.7544326019072023
$ws.Range("H16").Value = 1.394223664128001
$ws.Range("I16").Value = 1.105904788684794
$ws.Range("J16").Value = 0.6299189688320024
$ws.Range("K16").Value = 0.8804315356160015
$ws.Range("L16").Value = 0.9376918773248007
$ws.Range("M16").Value = 1.001729390830933

$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9919558868816298
$ws.Range("D17").Value = 0.9952238808813222
$ws.Range("E17").Value = 0.9951313150131709
$ws.Range("F17").Value = 0.9919558868816298
$ws.Range("G17").Value = 0.992536226054725
$ws.Range("H17").Value = 0.9965920879885471
$ws.Range("I17").Value = 0.9942953092046359
$ws.Range("J17").Value = 0.9952238808813222
$ws.Range("K17").Value = 0.9951775979472466
$ws.Range("L17").Value = 0.9935667424144382
$ws.Range("M17").Value = 0.994289117670672

$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9948006431056027
$ws.Range("D18").Value = 1.016072100310265
$ws.Range("E18").Value = 0.9847001417652971
$ws.Range("F18").Value = 0.9948006431056027
$ws.Range("G18").Value = 1.004898962236324
$ws.Range("H18").Value = 0.9721639649292282
$ws.Range("I18").Value = 0.9895120078572394
$ws.Range("J18").Value = 1.016072100310265
$ws.Range("K18").Value = 1.000386121037781
$ws.Range("L18").Value = 0.9975933820716918
$ws.Range("M18").Value = 0.993691303367326

$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9801333066439487
$ws.Range("D19").Value = 1.071366385799256
$ws.Range("E19").Value = 0.9739668093533926
$ws.Range("F19").Value = 0.9801333066439487
$ws.Range("G19").Value = 1.038761087289167
$ws.Range("H19").Value = 0.9339298383575777
$ws.Range("I19").Value = 0.9726723478379107
$ws.Range("J19").Value = 1.071366385799256
$ws.Range("K19").Value = 1.022666597576324
$ws.Range("L19").Value = 1.001399952110136
$ws.Range("M19").Value = 0.995138295880209

# New rows 17-19: set index column and copy header/index cell style from row 16
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
